$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J12").Value = "vfvdfvdvdfVdvdvd"
$ws.Range("J13").Value = "vfvdfvdvdfVdvdvd"
$ws.Range("J14").Value = "dv"
$ws.Range("J15").Value = "t"
$ws.Range("J16").Value = "vrt"
$ws.Range("J17").Value = "rt"
$ws.Range("J18").Value = "rt"
$ws.Range("J19").Value = "t"

$null = $ws.Range("J20").Select()
